$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.780.35'
$ws.Range("E2").Value = '  +4.52%  '
$ws.Range("D3").Value = '2.730.81'
$ws.Range("E3").Value = '  +2.68%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''577.41'
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("D6").Value = '''154.20'
$ws.Range("E6").Value = '  +5.93%  '
$ws.Range("D7").Value = '''0.995'
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("D8").Value = '''0.611'
$ws.Range("E8").Value = '  +1.74%  '
$ws.Range("D9").Value = '2.756.54'
$ws.Range("D10").Value = '''6.70'
$ws.Range("E11").Value = '  +4.52%  '
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").Value = '''0.388'
$ws.Range("E12").Value = '  +2.02%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '''0.161'
$ws.Range("E13").Value = '  +4.49%  '
$ws.Range("D14").Value = '3.216.55'
$ws.Range("E14").Value = '  +2.73%  '
$ws.Range("D15").Value = '''26.29'
$ws.Range("E15").Value = '  +2.46%  '
$ws.Range("D16").Value = '63.664.01'
$ws.Range("E16").Value = '  +4.39%  '
$ws.Range("E17").Value = '  +5.52%  '
$ws.Range("D18").Value = '2.748.99'
$ws.Range("E18").Value = '  +2.74%  '
$ws.Range("D19").Value = '''11.95'
$ws.Range("E19").Value = '  +2.64%  '
$ws.Range("E20").Value = '  +2.41%  '
$ws.Range("D21").Value = '''360.61'
$ws.Range("E21").Value = '  +2.32%  '
$ws.Range("D22").Value = '''6.95'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").Value = '''0.997'
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("E24").Value = '  -0.26%  '
$ws.Range("D25").Value = '''66.09'
$ws.Range("E25").Value = '  +2.99%  '
$ws.Range("D26").Value = '''0.169'
$ws.Range("E26").Value = '  +4.71%  '
$ws.Range("D27").Value = '''8.50'
$ws.Range("E27").Value = '  +3.36%  '
$ws.Range("D28").Value = '''0.997'
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("D29").Value = '0.0₃0912'
$ws.Range("E29").Value = '  +11.00%  '
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("D31").Value = '''7.10'
$ws.Range("E31").Value = '  +4.44%  '
$ws.Range("E32").Value = '  +2.45%  '
$ws.Range("E33").Value = '  +14.10%  '
$ws.Range("D34").Value = '''0.996'
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").Value = '''20.45'
$ws.Range("E35").Value = '  +2.55%  '
$ws.Range("D36").Value = '''4.77'
$ws.Range("E36").Value = '  +6.92%  '
$ws.Range("E37").Value = '  +8.32%  '
$ws.Range("E38").Value = '  +8.87%  '
$ws.Range("E39").Value = '  +14.47%  '
$ws.Range("D40").Value = '''345.42'
$ws.Range("E40").Value = '  +5.43%  '
$ws.Range("E41").Value = '  +5.08%  '
$ws.Range("D42").Value = '''39.38'
$ws.Range("E42").Value = '  +2.23%  '
$ws.Range("D43").Value = '''5.61'
$ws.Range("E43").Value = '  +6.68%  '
$ws.Range("D44").Value = '''21.72'
$ws.Range("E44").Value = '  +5.16%  '
$ws.Range("D45").Value = '''21.71'
$ws.Range("E45").Value = '  +5.76%  '
$ws.Range("E46").Value = '  +5.19%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '''0.646'
$ws.Range("E47").Value = '  +4.76%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '''139.05'
$ws.Range("E48").Value = '  +3.09%  '
$ws.Range("D49").Value = '''0.0254'
$ws.Range("E49").Value = '  +2.86%  '
$ws.Range("E50").Value = '  +0.55%  '
$ws.Range("D51").Value = '''0.995'
$ws.Range("E51").Value = '  -0.25%  '
